$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column B (old B shifts to C)
$ws.Columns("B:B").Insert()

# Fill the new column B with ":" for data rows 2-47 (header B1 stays blank)
$range = $ws.Range("B2:B47")
$range.Value = ":"

# Set width of new column B
$ws.Columns("B:B").ColumnWidth = 5.42578125

# Update the autofilter range to cover the new column
$ws.Range("A1:C47").AutoFilter(1)

# Add a new empty sheet named Sheet2
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Sheet2"

# Move selection / view
$ws.Range("A2:C47").Select()
